$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.162.06"
Set-TextValue $ws.Range("E2") "  -0.32%  "
Set-TextValue $ws.Range("D3") "1.833.05"
Set-TextValue $ws.Range("E3") "  +1.59%  "
Set-TextValue $ws.Range("D4") "0.9985"
Set-TextValue $ws.Range("E4") "  -0.54%  "
Set-TextValue $ws.Range("D5") "311.36"
Set-TextValue $ws.Range("E5") "  -0.90%  "
Set-TextValue $ws.Range("D6") "0.9995"
Set-TextValue $ws.Range("E6") "  -0.31%  "
Set-TextValue $ws.Range("D7") "0.5129"
Set-TextValue $ws.Range("E7") "  -2.78%  "
Set-TextValue $ws.Range("D8") "0.3950"
Set-TextValue $ws.Range("E8") "  +3.19%  "
Set-TextValue $ws.Range("D9") "0.09816"
Set-TextValue $ws.Range("E9") "  +22.50%  "
Set-TextValue $ws.Range("D10") "1.109"
Set-TextValue $ws.Range("E10") "  +0.80%  "
Set-TextValue $ws.Range("D11") "40.83"
Set-TextValue $ws.Range("E11") "  -1.42%  "
Set-TextValue $ws.Range("D12") "6.475"
Set-TextValue $ws.Range("E12") "  +2.35%  "
Set-TextValue $ws.Range("D13") "0.9991"
Set-TextValue $ws.Range("E13") "  -0.43%  "
Set-TextValue $ws.Range("D14") "20.58"
Set-TextValue $ws.Range("E14") "  -0.08%  "
Set-TextValue $ws.Range("B15") "Chainlink"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D15") "7.413"
Set-TextValue $ws.Range("E15") "  +1.25%  "
Set-TextValue $ws.Range("B16") "WrappedEther"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "1.809.29"
Set-TextValue $ws.Range("E16") "  +0.06%  "
Set-TextValue $ws.Range("D17") "0.00001143"
Set-TextValue $ws.Range("E17") "  +4.26%  "
Set-TextValue $ws.Range("D18") "93.17"
Set-TextValue $ws.Range("E18") "  +1.10%  "
Set-TextValue $ws.Range("D19") "0.06602"
Set-TextValue $ws.Range("E19") "  -0.12%  "
Set-TextValue $ws.Range("B20") "Avalanche"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D20") "17.40"
Set-TextValue $ws.Range("E20") "  +0.06%  "
Set-TextValue $ws.Range("B21") "Dai"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D21") "0.9995"
Set-TextValue $ws.Range("E21") "  -0.34%  "
Set-TextValue $ws.Range("D22") "6.066"
Set-TextValue $ws.Range("E22") "  +1.69%  "
Set-TextValue $ws.Range("D23") "28.218.91"
Set-TextValue $ws.Range("E23") "  -0.33%  "
Set-TextValue $ws.Range("D24") "11.17"
Set-TextValue $ws.Range("E24") "  +0.11%  "
Set-TextValue $ws.Range("D25") "2.249"
Set-TextValue $ws.Range("E25") "  -0.53%  "
Set-TextValue $ws.Range("B26") "EthereumClassic"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D26") "20.73"
Set-TextValue $ws.Range("E26") "  +1.23%  "
Set-TextValue $ws.Range("B27") "LidoDAOToken"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "2.448"
Set-TextValue $ws.Range("E27") "  +3.78%  "
Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "157.24"
Set-TextValue $ws.Range("E28") "  -2.14%  "
Set-TextValue $ws.Range("B29") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D29") "2.029.90"
Set-TextValue $ws.Range("E29") "  +0.90%  "
Set-TextValue $ws.Range("D30") "129.12"
Set-TextValue $ws.Range("E30") "  +4.64%  "
Set-TextValue $ws.Range("D31") "0.1096"
Set-TextValue $ws.Range("E31") "  +0.89%  "
Set-TextValue $ws.Range("D32") "1.060"
Set-TextValue $ws.Range("E32") "  +0.34%  "
Set-TextValue $ws.Range("D33") "5.655"
Set-TextValue $ws.Range("E33") "  +1.76%  "
Set-TextValue $ws.Range("D34") "3.632"
Set-TextValue $ws.Range("E34") "  -1.52%  "
Set-TextValue $ws.Range("D35") "0.06932"
Set-TextValue $ws.Range("E35") "  -4.41%  "
Set-TextValue $ws.Range("D36") "9.100"
Set-TextValue $ws.Range("E36") "  +5.04%  "
Set-TextValue $ws.Range("D37") "0.02347"
Set-TextValue $ws.Range("E37") "  +1.21%  "
Set-TextValue $ws.Range("D38") "0.2182"
Set-TextValue $ws.Range("E38") "  +1.15%  "
Set-TextValue $ws.Range("D39") "11.61"
Set-TextValue $ws.Range("E39") "  -6.08%  "
Set-TextValue $ws.Range("D40") "5.024"
Set-TextValue $ws.Range("E40") "  -1.80%  "
Set-TextValue $ws.Range("D41") "0.6262"
Set-TextValue $ws.Range("E41") "  +0.90%  "
Set-TextValue $ws.Range("D42") "0.9990"
Set-TextValue $ws.Range("D43") "1.154"
Set-TextValue $ws.Range("E43") "  -1.14%  "
Set-TextValue $ws.Range("D44") "13.36"
Set-TextValue $ws.Range("E44") "  +0.31%  "
Set-TextValue $ws.Range("D45") "0.5995"
Set-TextValue $ws.Range("E45") "  -0.42%  "
Set-TextValue $ws.Range("E46") "  -5.83%  "
Set-TextValue $ws.Range("D47") "3.711"
Set-TextValue $ws.Range("E47") "  -1.47%  "
Set-TextValue $ws.Range("D48") "125.47"
Set-TextValue $ws.Range("E48") "  -1.29%  "
Set-TextValue $ws.Range("D49") "1.977"
Set-TextValue $ws.Range("E49") "  +2.46%  "
Set-TextValue $ws.Range("D50") "1.189"
Set-TextValue $ws.Range("E50") "  -2.30%  "
Set-TextValue $ws.Range("D51") "0.06788"
Set-TextValue $ws.Range("E51") "  -0.47%  "
